$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"
$metaPara.Range.Text = "Meta description: Read our review of CherryPop online slot game. Play CherryPop for free or real money and learn about RTP, bonuses, and gameplay features."

# Bold just the "Meta description" label (first 16 characters).
$metaStart = $d.Paragraphs(2).Range.Start
$labelRange = $d.Range($metaStart, $metaStart + 16)
$labelRange.Bold = 1

# ---------------------------------------------------------------------------
# 2) Near the end of the document: drop the duplicated bold title paragraph,
#    and rewrite the italic "meta description" paragraph into an image
#    prompt paragraph (keeping the italic formatting).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Play CherryPop Slot for Free - Review & Real Money RTP" -and $i -ne 1) {
        $p.Range.Delete()
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Read our review of CherryPop online slot game*") {
        $start = $p.Range.Start
        $end = $p.Range.End
        $r = $d.Range($start, $end)
        $r.Text = "Create a feature image that captures the essence of CherryPop, a fun and exciting online slot game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior can be holding a cherry or surrounded by exploding symbols to represent the PopWins mode. The background should be colorful and vibrant, with elements of fruits, gems, and playing cards to represent the different symbols in the game. The image should also have the CherryPop logo prominently displayed, along with the tagline `"Pop your way to massive winnings!`""
        break
    }
}
